$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '28.466.84'
$ws.Range("E2").Value = '  +2.06%  '
$ws.Range("D3").Value = '1.575.66'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  +1.08%  '
Set-TextValue "D5" '210.97'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E7").Value = '  +1.00%  '
Set-TextValue "D8" '46.22'
$ws.Range("E8").Value = '  +4.68%  '
Set-TextValue "D9" '23.80'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '1.799.47'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '1.565.59'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("E15").Value = '  +0.48%  '
Set-TextValue "D16" '3.71'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '28.449.01'
$ws.Range("E17").Value = '  +2.00%  '
Set-TextValue "D18" '62.36'
$ws.Range("E18").Value = '  -1.67%  '
Set-TextValue "D19" '229.20'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").Value = '0.0₃0693'
$ws.Range("E22").Value = '  +1.10%  '
Set-TextValue "D23" '3.93'
$ws.Range("E23").Value = '  -4.38%  '
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("E25").Value = '  +3.84%  '
Set-TextValue "D26" '150.75'
$ws.Range("E26").Value = '  -0.25%  '
Set-TextValue "D27" '15.02'
$ws.Range("E27").Value = '  -1.40%  '
Set-TextValue "D28" '6.46'
$ws.Range("E28").Value = '  -1.54%  '
Set-TextValue "D29" '0.105'
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("E30").Value = '  +1.04%  '
Set-TextValue "D31" '1.11'
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("E33").Value = '  -0.65%  '
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("D35").Value = '1.393.34'
$ws.Range("E35").Value = '  -1.76%  '
$ws.Range("E36").Value = '  -2.00%  '
Set-TextValue "D37" '1.01'
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("E38").Value = '  +3.15%  '
Set-TextValue "D39" '0.0166'
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").Value = '  +3.89%  '
Set-TextValue "D41" '0.532'
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("E42").Value = '  +1.17%  '
Set-TextValue "D43" '0.793'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("E44").Value = '  +0.31%  '
Set-TextValue "D45" '1.85'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("E46").Value = '  +1.22%  '
Set-TextValue "D47" '62.38'
$ws.Range("E47").Value = '  -2.37%  '
$ws.Range("D48").Value = '1.711.50'
$ws.Range("E48").Value = '  +0.43%  '
Set-TextValue "D49" '85.77'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("E50").Value = '  +4.52%  '
Set-TextValue "D51" '0.0517'
$ws.Range("E51").Value = '  -1.22%  '
